# Update the "Price" (D) and "Volume(1h)" (E) columns with the latest scraped
# values, and fix the Mantle/Aptos row ordering (rows 50-51 swapped places).
# NOTE: several Price values (e.g. "0.5586", "1.008") look like plain numbers,
# so a leading apostrophe is used to force Excel to store them as text,
# exactly like the original data (which is text, not numeric).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '26.366.43'
$ws.Range("E2").Value = '  +1.24%  '
$ws.Range("D3").Value = '1.685.37'
$ws.Range("E3").Value = '  +1.04%  '
$ws.Range("D4").Value = '''1.008'
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("E5").Value = '  +0.78%  '
$ws.Range("D6").Value = '''0.5586'
$ws.Range("E6").Value = '  +9.64%  '
$ws.Range("D8").Value = '''0.2711'
$ws.Range("E8").Value = '  +1.93%  '
$ws.Range("D9").Value = '''0.06508'
$ws.Range("E9").Value = '  +1.73%  '
$ws.Range("D10").Value = '''22.16'
$ws.Range("E10").Value = '  +1.83%  '
$ws.Range("D11").Value = '''0.07568'
$ws.Range("E11").Value = '  +1.67%  '
$ws.Range("D12").Value = '''4.553'
$ws.Range("E12").Value = '  +0.86%  '
$ws.Range("D13").Value = '1.683.00'
$ws.Range("E13").Value = '  +0.60%  '
$ws.Range("D14").Value = '''0.5822'
$ws.Range("E14").Value = '  -0.16%  '
$ws.Range("D15").Value = '''0.000008477'
$ws.Range("E15").Value = '  -0.75%  '
$ws.Range("D16").Value = '''65.44'
$ws.Range("E16").Value = '  +1.78%  '
$ws.Range("D17").Value = '26.400.05'
$ws.Range("E17").Value = '  +1.14%  '
$ws.Range("D18").Value = '''4.952'
$ws.Range("E18").Value = '  +0.24%  '
$ws.Range("E19").Value = '  +0.24%  '
$ws.Range("E20").Value = '  +1.43%  '
$ws.Range("D21").Value = '''191.66'
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("D22").Value = '''6.249'
$ws.Range("E22").Value = '  +0.83%  '
$ws.Range("D23").Value = '''1.008'
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("D24").Value = '''148.67'
$ws.Range("E24").Value = '  +2.76%  '
$ws.Range("D25").Value = '''0.1334'
$ws.Range("E25").Value = '  +11.38%  '
$ws.Range("D26").Value = '''7.934'
$ws.Range("E26").Value = '  +4.25%  '
$ws.Range("D27").Value = '''15.86'
$ws.Range("E27").Value = '  +1.34%  '
$ws.Range("D28").Value = '''0.06327'
$ws.Range("E28").Value = '  -3.45%  '
$ws.Range("D29").Value = '''1.393'
$ws.Range("E29").Value = '  +4.07%  '
$ws.Range("D30").Value = '''1.324'
$ws.Range("E30").Value = '  +0.69%  '
$ws.Range("D31").Value = '''3.600'
$ws.Range("E31").Value = '  +1.68%  '
$ws.Range("E32").Value = '  +1.95%  '
$ws.Range("E33").Value = '  +1.24%  '
$ws.Range("E34").Value = '  +2.38%  '
$ws.Range("D35").Value = '''0.6232'
$ws.Range("E35").Value = '  +1.83%  '
$ws.Range("D36").Value = '''2.399'
$ws.Range("E36").Value = '  +1.24%  '
$ws.Range("D37").Value = '''2.719'
$ws.Range("E37").Value = '  +1.29%  '
$ws.Range("D38").Value = '''6.235'
$ws.Range("E38").Value = '  -0.57%  '
$ws.Range("D39").Value = '1.115.81'
$ws.Range("E39").Value = '  +2.25%  '
$ws.Range("D40").Value = '''0.01633'
$ws.Range("E40").Value = '  +2.14%  '
$ws.Range("E41").Value = '  +0.88%  '
$ws.Range("E42").Value = '  +0.44%  '
$ws.Range("D43").Value = '''100.68'
$ws.Range("E43").Value = '  -0.49%  '
$ws.Range("D44").Value = '1.834.55'
$ws.Range("E44").Value = '  +1.01%  '
$ws.Range("E45").Value = '  -4.89%  '
$ws.Range("D46").Value = '''57.45'
$ws.Range("E46").Value = '  +1.86%  '
$ws.Range("D47").Value = '''8.191'
$ws.Range("E47").Value = '  +1.43%  '
$ws.Range("E48").Value = '  -0.06%  '
$ws.Range("E49").Value = '  +0.80%  '
$ws.Range("B50").Value = 'Aptos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D50").Value = '''6.101'
$ws.Range("E50").Value = '  +0.74%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = '''0.4297'
$ws.Range("E51").Value = '  +0.23%  '
